# edit.ps1 — Append 50 new "Filename" labels to the data-thesis tracking sheet
# (rows 52-101), matching the commit "Added 50 labels; Total 100".
#
# Column layout (row 1 header): A=Filename, B..I = algorithm flags
# (Quicksort, Mergesort, Selectionsort, Insertionsort, Bubblesort,
# Linear search, Binary Search, Linked List), J = Hashmap flag.
# New rows keep B..I at 0 (not reviewed for those algorithms yet) and set
# J to 0/1 per whether the file implements a hashmap.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    "DirectoryImpl.java",
    "DivisionMethod.java",
    "double_hash (2).py",
    "double_hash (3).py",
    "double_hash (4).py",
    "double_hash.py",
    "double_hashing.py",
    "DoubleHashing.java",
    "DoublyLinkedList.java",
    "Driver.java",
    "dynamic_hashing.py",
    "Entry.java",
    "EntryIterator.java",
    "EntrySet.java",
    "FastMap.java",
    "first-recurring-number.js",
    "FixedSizedHashMap.java",
    "hash_table (2).py",
    "hash_table (3).py",
    "hash_table (4).py",
    "hash_table (5).py",
    "hash_table_with_linked_list (2).py",
    "hash_table_with_linked_list (3).py",
    "hash_table_with_linked_list (4).py",
    "hash_table_with_linked_list.py",
    "hash_table.py",
    "hash-table.js",
    "Hash.js",
    "HashArrayMappedTrie.java",
    "HashCodeAndEquals.java",
    "HashEntry.java",
    "HashFunction (2).java",
    "HashFunction.java",
    "HashMap (2).java",
    "HashMap (3).java",
    "HashMap (4).java",
    "HashMap (5).java",
    "HashMap (6).java",
    "HashMap (7).java",
    "HashMap (8).java",
    "HashMap (9).java",
    "HashMap (10).java",
    "HashMap (11).java",
    "HashMap (12).java",
    "HashMap (13).java",
    "HashMap (14).java",
    "HashMap.java",
    "HashMap.js",
    "hashmap.py",
    "HashMapDriver.java"
)

$jvals = @(0,0,1,1,1,1,1,0,1,1,1,0,1,1,1,0,1,1,1,1,1,1,1,1,1,1,1,0,1,1,1,0,0,0,1,0,1,1,1,1,0,1,1,1,1,0,0,1,1,1)

$startRow = 52
for ($i = 0; $i -lt $names.Count; $i++) {
    $row = $startRow + $i

    $ws.Cells.Item($row, 1).Value = $names[$i]

    for ($col = 2; $col -le 9; $col++) {
        $ws.Cells.Item($row, $col).Value = 0
    }

    $ws.Cells.Item($row, 10).Value = $jvals[$i]
}

$lastRow = $startRow + $names.Count - 1
$ws.Range("J$lastRow").Select() | Out-Null
